$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style-template cells already present in the sheet (column B), one per
# style index used by the new column D so PasteSpecial(xlPasteFormats)
# reuses the existing style index instead of minting new ones.
$tmpl = @{1='B5'; 2='B2'; 3='B3'}

# Row-by-row data for the new column D (mirrors column C, refreshed at a
# later timestamp -- some values identical, some updated).
$rows = @(
    @{r=1; style=1; type='s'; val='2025/10/22'},
    @{r=2; style=2; type='s'; val='上证'},
    @{r=3; style=3; type='n'; val=63.83},
    @{r=4; style=3; type='n'; val=3914.49},
    @{r=5; style=1; type='e'; val=$null},
    @{r=6; style=3; type='n'; val=50.06},
    @{r=7; style=3; type='n'; val=5498.71},
    @{r=8; style=1; type='e'; val=$null},
    @{r=9; style=3; type='n'; val=55.75},
    @{r=10; style=3; type='n'; val=4596.44},
    @{r=11; style=1; type='e'; val=$null},
    @{r=12; style=3; type='n'; val=61.67},
    @{r=13; style=3; type='n'; val=7131.91},
    @{r=14; style=1; type='e'; val=$null},
    @{r=15; style=3; type='n'; val=28.36},
    @{r=16; style=3; type='n'; val=2675.88},
    @{r=17; style=1; type='e'; val=$null},
    @{r=18; style=3; type='n'; val=95.76000000000001},
    @{r=19; style=3; type='n'; val=6735.35},
    @{r=20; style=1; type='e'; val=$null},
    @{r=21; style=3; type='n'; val=67.68000000000001},
    @{r=22; style=3; type='n'; val=84426.34},
    @{r=23; style=1; type='e'; val=$null},
    @{r=24; style=3; type='n'; val=85.68000000000001},
    @{r=25; style=3; type='n'; val=19909.14},
    @{r=26; style=1; type='e'; val=$null},
    @{r=27; style=3; type='n'; val=78.34999999999999},
    @{r=28; style=3; type='n'; val=39894.54},
    @{r=29; style=1; type='e'; val=$null},
    @{r=30; style=3; type='n'; val=57.37},
    @{r=31; style=3; type='n'; val=5649.26},
    @{r=32; style=1; type='e'; val=$null},
    @{r=33; style=3; type='n'; val=10.04},
    @{r=34; style=3; type='n'; val=33455.38},
    @{r=35; style=1; type='e'; val=$null},
    @{r=36; style=3; type='n'; val=30.5},
    @{r=37; style=3; type='n'; val=3229.07},
    @{r=38; style=1; type='e'; val=$null},
    @{r=39; style=3; type='n'; val=48.4},
    @{r=40; style=3; type='n'; val=3059.91},
    @{r=41; style=1; type='e'; val=$null},
    @{r=42; style=3; type='n'; val=19.53},
    @{r=43; style=3; type='n'; val=7397.02},
    @{r=44; style=1; type='e'; val=$null},
    @{r=45; style=3; type='n'; val=33.2},
    @{r=46; style=3; type='n'; val=9016.719999999999},
    @{r=47; style=1; type='e'; val=$null},
    @{r=48; style=3; type='n'; val=9.710000000000001},
    @{r=49; style=3; type='n'; val=12976.73},
    @{r=50; style=1; type='e'; val=$null},
    @{r=51; style=3; type='n'; val=22.73},
    @{r=52; style=3; type='n'; val=12435.27},
    @{r=53; style=1; type='e'; val=$null},
    @{r=54; style=3; type='n'; val=18.98},
    @{r=55; style=3; type='n'; val=9752.719999999999},
    @{r=56; style=1; type='e'; val=$null},
    @{r=57; style=3; type='n'; val=24.91},
    @{r=58; style=3; type='n'; val=16165.24},
    @{r=59; style=1; type='e'; val=$null},
    @{r=60; style=3; type='n'; val=33.53},
    @{r=61; style=3; type='n'; val=17526.85},
    @{r=62; style=1; type='e'; val=$null},
    @{r=63; style=3; type='n'; val=21.28},
    @{r=64; style=3; type='n'; val=10228.42},
    @{r=65; style=1; type='e'; val=$null},
    @{r=66; style=3; type='n'; val=15.36},
    @{r=67; style=3; type='n'; val=9930.68},
    @{r=68; style=1; type='e'; val=$null},
    @{r=69; style=3; type='n'; val=21.08},
    @{r=70; style=3; type='n'; val=3142.01},
    @{r=71; style=1; type='e'; val=$null},
    @{r=72; style=3; type='n'; val=44.37},
    @{r=73; style=3; type='n'; val=5926.91},
    @{r=74; style=1; type='e'; val=$null},
    @{r=75; style=3; type='n'; val=30.39},
    @{r=76; style=3; type='n'; val=9486.709999999999},
    @{r=77; style=1; type='e'; val=$null},
    @{r=78; style=3; type='n'; val=18.35},
    @{r=79; style=3; type='n'; val=2458.76},
    @{r=80; style=1; type='e'; val=$null},
    @{r=81; style=3; type='n'; val=55.78},
    @{r=82; style=3; type='n'; val=2635.51},
    @{r=83; style=1; type='e'; val=$null},
    @{r=84; style=3; type='n'; val=58.79},
    @{r=85; style=3; type='n'; val=2698.2},
    @{r=86; style=1; type='e'; val=$null},
    @{r=87; style=3; type='n'; val=51.71},
    @{r=88; style=3; type='n'; val=3783.74},
    @{r=89; style=1; type='e'; val=$null},
    @{r=90; style=3; type='n'; val=48.63},
    @{r=91; style=3; type='n'; val=2067.64},
    @{r=92; style=1; type='e'; val=$null},
    @{r=93; style=3; type='n'; val=27.72},
    @{r=94; style=3; type='n'; val=13758.07},
    @{r=95; style=1; type='e'; val=$null},
    @{r=96; style=3; type='n'; val=88.04000000000001},
    @{r=97; style=3; type='n'; val=9221.57},
    @{r=98; style=1; type='e'; val=$null},
    @{r=99; style=3; type='n'; val=58.2},
    @{r=100; style=3; type='n'; val=12069.97},
    @{r=101; style=1; type='e'; val=$null},
    @{r=102; style=3; type='n'; val=3.09},
    @{r=103; style=3; type='n'; val=2235.14},
    @{r=104; style=1; type='e'; val=$null},
    @{r=105; style=3; type='n'; val=30.47},
    @{r=106; style=3; type='n'; val=869.04},
    @{r=107; style=1; type='e'; val=$null},
    @{r=108; style=3; type='n'; val=29.71},
    @{r=109; style=3; type='n'; val=2659.06},
    @{r=110; style=1; type='e'; val=$null},
    @{r=111; style=3; type='n'; val=20.62},
    @{r=112; style=3; type='n'; val=4050.56},
    @{r=113; style=1; type='e'; val=$null},
    @{r=114; style=3; type='n'; val=29.02},
    @{r=115; style=3; type='n'; val=3503.12}
)

foreach ($row in $rows) {
    $ref = "D" + $row.r
    $t = $tmpl[$row.style]

    if ($row.type -eq 's') {
        # Text value (date label / header) -- force Text format first so
        # Excel doesn't auto-convert a date-looking string into a serial
        # date number.
        $ws.Range($ref).NumberFormat = "@"
        $ws.Range($ref).Value = $row.val
    } elseif ($row.type -eq 'n') {
        $ws.Range($ref).Value = $row.val
    }
    # else 'e' -> leave the cell empty, only formatting is applied below.

    $ws.Range($t).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
